$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column header corrections (row 1)
$ws.Range("A1").Value = "gender"
$ws.Range("B1").Value = "age"
$ws.Range("C1").Value = "university year"
$ws.Range("D1").Value = "faculty"
$ws.Range("E1").Value = "asm"
$ws.Range("F1").Value = "aum"
$ws.Range("G1").Value = "cwc"

# Recode gender column (A2:A30): M -> 0, F -> 1
for ($r = 2; $r -le 30; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value() -eq "F") {
        $cell.Value = 1
    } else {
        $cell.Value = 0
    }
}
